# Apply updated forecast data: dates shift forward one week and several
# forecast figures are recalculated (see commit "Fixed update to excel issue").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Ensure the Week_Start_Date column keeps storing plain text (not auto-converted
# to a date serial number) by forcing a text number format before assigning.
$ws1.Range("B2:B17").NumberFormat = "@"

# --- Sheet "Forecast Comparison" ---

# Row 2 (W1)
$ws1.Range("B2").Value = "2025-02-02"
$ws1.Range("E2").Value = 3
$ws1.Range("F2").Value = 2
$ws1.Range("G2").Value = 5
$ws1.Range("H2").Value = 11

# Row 3 (W2)
$ws1.Range("B3").Value = "2025-02-09"

# Row 4 (W3)
$ws1.Range("B4").Value = "2025-02-16"

# Row 5 (W4)
$ws1.Range("B5").Value = "2025-02-23"
$ws1.Range("E5").Value = 5

# Row 6 (W5)
$ws1.Range("B6").Value = "2025-03-02"

# Row 7 (W6)
$ws1.Range("B7").Value = "2025-03-09"
$ws1.Range("H7").Value = 13

# Row 8 (W7)
$ws1.Range("B8").Value = "2025-03-16"
$ws1.Range("H8").Value = 13

# Row 9 (W8)
$ws1.Range("B9").Value = "2025-03-23"
$ws1.Range("G9").Value = 5
$ws1.Range("H9").Value = 12

# Row 10 (W9)
$ws1.Range("B10").Value = "2025-03-30"

# Row 11 (W10)
$ws1.Range("B11").Value = "2025-04-06"

# Row 12 (W11)
$ws1.Range("B12").Value = "2025-04-13"
$ws1.Range("D12").Value = 0
$ws1.Range("E12").Value = 3
$ws1.Range("F12").Value = 2

# Row 13 (W12)
$ws1.Range("B13").Value = "2025-04-20"

# Row 14 (W13)
$ws1.Range("B14").Value = "2025-04-27"

# Row 15 (W14)
$ws1.Range("B15").Value = "2025-05-04"
$ws1.Range("H15").Value = 9

# Row 16 (W15)
$ws1.Range("B16").Value = "2025-05-11"

# Row 17 (W16)
$ws1.Range("B17").Value = "2025-05-18"
$ws1.Range("G17").Value = 5

# --- Sheet "Summary" ---

$ws2.Range("B2:B15").NumberFormat = "@"

$ws2.Range("B2").Value = "2023-01-01 to 2025-01-26"
$ws2.Range("B5").Value = "22"
$ws2.Range("B7").Value = "32"
$ws2.Range("B9").Value = "8"
$ws2.Range("B13").Value = "2025-02-02"
$ws2.Range("B15").Value = "2025-03-30"

$wb.Save()
